# Update the "取得日時" (retrieved at) timestamp in column A for all data
# rows on the "ランサーズ" sheet to reflect the latest append run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-07 12:48:17"

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
